$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: mark C9 and E9 as correct (value 5) and highlight them the same way
# E11 already is (no fill / style index "2" look), matching the grading sheet's
# "un-greyed / reviewed" cell styling.
$ws.Range("E11").Copy() | Out-Null
$ws.Range("C9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C9").Value = 5
$ws.Range("E9").Value = 5

# Row 21: mark C21 and D21 as correct (value 5) with the same styling.
$ws.Range("E11").Copy() | Out-Null
$ws.Range("C21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 5

# Move the active selection to E21 (last cell touched while grading).
$ws.Range("E21").Select() | Out-Null

$excel.CutCopyMode = $false
